# add "Save" column in s_vals sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header cell (G1) onto the new
# header cell (H1), then overwrite its text with "Save".
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
